# Permute the per-record data (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Precio $/Kg, Kg/unidad) across rows 2..26.
# The mapping below gives, for each destination row, the row whose original
# values should be copied into it (derived from the target diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns that move together as one "record" block
$cols = @('D','L','M','N','O','P','Q','S','T')

# destination row -> source row (values read BEFORE any writes happen)
$mapping = @{
    2  = 11
    3  = 12
    4  = 16
    5  = 20
    6  = 7
    7  = 17
    8  = 18
    9  = 22
    10 = 3
    11 = 4
    12 = 9
    13 = 2
    14 = 19
    15 = 6
    16 = 24
    17 = 25
    18 = 10
    19 = 26
    20 = 15
    21 = 23
    22 = 5
    23 = 14
    24 = 8
    25 = 21
    26 = 13
}

# Snapshot the original values for every row/column before writing anything,
# since several destination rows read from rows that are also destinations.
$snapshot = @{}
foreach ($row in 2..26) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowVals
}

foreach ($destRow in 2..26) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
